$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$pane = $win.ActivePane
$pane | Get-Member
